$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01179113890272636
$ws.Range("C2").Value = 0.9527065656392815
$ws.Range("D2").Value = 4.106302077608952
$ws.Range("E2").Value = 2.026401262733754
$ws.Range("F2").Value = 2.046136859441977
$ws.Range("G2").Value = 52
$ws.Range("B3").Value = 0.06492620600918006
$ws.Range("C3").Value = 1.070932841330739
$ws.Range("D3").Value = 4.414872682735234
$ws.Range("E3").Value = 2.101159842262181
$ws.Range("F3").Value = 2.121054081328115
$ws.Range("G3").Value = 51
$ws.Range("B4").Value = 0.0188834649839998
$ws.Range("C4").Value = 0.9530525523764629
$ws.Range("D4").Value = 3.865729941648752
$ws.Range("E4").Value = 1.966145961430319
$ws.Range("F4").Value = 1.986015741865533
$ws.Range("G4").Value = 50
$ws.Range("B5").Value = 0.09692881375077184
$ws.Range("C5").Value = 1.09440480910644
$ws.Range("D5").Value = 4.629837775659457
$ws.Range("E5").Value = 2.151705782782455
$ws.Range("F5").Value = 2.171796906056384
$ws.Range("G5").Value = 49
$ws.Range("B6").Value = -0.000755637588625957
$ws.Range("C6").Value = 0.9503968899563854
$ws.Range("D6").Value = 3.951250381626263
$ws.Range("E6").Value = 1.987775234181738
$ws.Range("F6").Value = 2.0088103333927
$ws.Range("G6").Value = 48
$ws.Range("B7").Value = 0.03928990789677742
$ws.Range("C7").Value = 1.033046348395609
$ws.Range("D7").Value = 4.93222755344469
$ws.Range("E7").Value = 2.220861894275439
$ws.Range("F7").Value = 2.252012552851068
$ws.Range("G7").Value = 36
$ws.Range("B8").Value = 0.03868193857282701
$ws.Range("C8").Value = 1.016268278851836
$ws.Range("D8").Value = 5.049011224985318
$ws.Range("E8").Value = 2.247000495101262
$ws.Range("F8").Value = 2.279467318072585
$ws.Range("G8").Value = 35
$ws.Range("B9").Value = 0.02578329123745685
$ws.Range("C9").Value = 1.456619164208787
$ws.Range("D9").Value = 9.170609166841203
$ws.Range("E9").Value = 3.028301366581801
$ws.Range("F9").Value = 3.115983453452396
$ws.Range("G9").Value = 18
$ws.Range("B10").Value = -0.8477958937427421
$ws.Range("C10").Value = 1.04746007141034
$ws.Range("D10").Value = 6.235465529401324
$ws.Range("E10").Value = 2.497091413905651
$ws.Range("F10").Value = 2.463407886881444
$ws.Range("G10").Value = 11
$ws.Range("B11").Value = -0.248566666372188
$ws.Range("C11").Value = 0.301374533835348
$ws.Range("D11").Value = 0.1561541005302556
$ws.Range("E11").Value = 0.3951633846021865
$ws.Range("F11").Value = 0.3434543508584382
$ws.Range("G11").Value = 5
